$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:5PAISA"
$ws.Range("C2").Value = "NSE:ALKALI"
$ws.Range("D2").Value = "NSE:BALRAMCHIN"
$ws.Range("F2").Value = "NSE:ASIANPAINT"
$ws.Range("B3").Value = "NSE:AARVEEDEN"
$ws.Range("C3").Value = "NSE:ALLSEC"
$ws.Range("D3").Value = "NSE:HDFCLIFE"
$ws.Range("F3").Value = "NSE:DIXON"
$ws.Range("B4").Value = "NSE:ADANIENT"
$ws.Range("C4").Value = "NSE:CANFINHOME"
$ws.Range("D4").Value = "NSE:HINDCOPPER"
$ws.Range("F4").Value = "NSE:GODREJPROP"
$ws.Range("B5").Value = "NSE:ADANIPOWER"
$ws.Range("C5").Value = "NSE:DHANI"
$ws.Range("D5").Value = "NSE:JSWSTEEL"
$ws.Range("F5").Value = "NSE:GRANULES"
$ws.Range("B6").Value = "NSE:ADL"
$ws.Range("C6").Value = "NSE:INDIACEM"
$ws.Range("D6").Value = "NSE:MARUTI"
$ws.Range("F6").Value = "NSE:MARUTI"
$ws.Range("B7").Value = "NSE:AJANTPHARM"
$ws.Range("C7").Value = "NSE:INTELLECT"
$ws.Range("F7").Value = "NSE:POLYCAB"
$ws.Range("B8").Value = "NSE:ANGELONE"
$ws.Range("C8").Value = "NSE:LANDMARK"
$ws.Range("B9").Value = "NSE:ARROWGREEN"
$ws.Range("C9").Value = "NSE:MAHLOG"
$ws.Range("B10").Value = "NSE:ASIANPAINT"
$ws.Range("C10").Value = "NSE:MAITHANALL"
$ws.Range("B11").Value = "NSE:AUTOIND"
$ws.Range("C11").Value = "NSE:MALLCOM"
$ws.Range("B12").Value = "NSE:AYMSYNTEX"
$ws.Range("C12").Value = "NSE:MHRIL"
$ws.Range("B13").Value = "NSE:BALPHARMA"
$ws.Range("C13").Value = "NSE:MTNL"
$ws.Range("B14").Value = "NSE:BSL"
$ws.Range("C14").Value = "NSE:ORICONENT"
$ws.Range("B15").Value = "NSE:CHAMBLFERT"
$ws.Range("C15").Value = "NSE:PGHH"
$ws.Range("B16").Value = "NSE:CRAFTSMAN"
$ws.Range("C16").Value = "NSE:PLAZACABLE"
$ws.Range("B17").Value = "NSE:DECCANCE"
$ws.Range("C17").Value = "NSE:RATNAVEER"
$ws.Range("B18").Value = "NSE:DIXON"
$ws.Range("B19").Value = "NSE:EDELWEISS"
$ws.Range("B20").Value = "NSE:FIBERWEB"
$ws.Range("B21").Value = "NSE:GAIL"
$ws.Range("B22").Value = "NSE:GRANULES"
$ws.Range("B23").Value = "NSE:GSFC"
$ws.Range("B24").Value = "NSE:HMVL"
$ws.Range("B25").Value = "NSE:HUBTOWN"
$ws.Range("B26").Value = "NSE:INSPIRISYS"
$ws.Range("B27").Value = "NSE:KANORICHEM"
$ws.Range("B28").Value = "NSE:KERNEX"
$ws.Range("B29").Value = "NSE:MARUTI"
$ws.Range("B30").Value = "NSE:MCL"
$ws.Range("B31").Value = "NSE:POLYCAB"
$ws.Range("B32").Value = "NSE:PREMIERPOL"
$ws.Range("B33").Value = "NSE:RAIN"
$ws.Rows("34:49").Delete()
